$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 07:22"

# Refresh the full country table (rows 4-212): COVID-19 case counts as of 07:22,
# re-sorted by total cases (column B) descending, matching the source feed order.
$data = @(
    @("Estados Unidos", 311635, 278, 14825, 288356, 8206, 2, 8454),
    @("España", 126168, 0, 34219, 80002, 6532, 0, 11947),
    @("Italia", 124632, 0, 20996, 88274, 3994, 0, 15362),
    @("Alemania", 96092, 0, 26400, 68248, 3936, 0, 1444),
    @("Francia", 89953, 0, 15438, 66955, 6838, 0, 7560),
    @("China", 81669, 30, 76964, 1376, 295, 3, 3329),
    @("Iran", 55743, 0, 19736, 32555, 4103, 0, 3452),
    @("Reino Unido", 41903, 0, 135, 37455, 163, 0, 4313),
    @("Turquia", 23934, 0, 786, 22647, 1311, 0, 501),
    @("Suiza", 20505, 0, 6415, 13424, 391, 0, 666),
    @("Belgica", 18431, 0, 3247, 13901, 1245, 0, 1283),
    @("Paises Bajos", 16627, 0, 250, 14726, 1360, 0, 1651),
    @("Canada", 14018, 106, 2603, 11182, 426, 2, 233),
    @("Austria", 11781, 0, 2507, 9088, 245, 0, 186),
    @("Portugal", 10524, 0, 75, 10183, 251, 0, 266),
    @("Brasil", 10360, 0, 127, 9788, 296, 0, 445),
    @("Corea del Sur", 10237, 81, 6463, 3591, 55, 6, 183),
    @("Israel", 7851, 0, 427, 7380, 115, 0, 44),
    @("Suecia", 6443, 0, 205, 5865, 379, 0, 373),
    @("Australia", 5635, 85, 585, 5016, 85, 4, 34),
    @("Noruega", 5550, 0, 32, 5456, 98, 0, 62),
    @("Rusia", 4731, 0, 333, 4355, 8, 0, 43),
    @("Irlanda", 4604, 0, 25, 4442, 148, 0, 137),
    @("Chequia", 4472, 0, 78, 4335, 87, 0, 59),
    @("Chile", 4161, 0, 528, 3606, 38, 0, 27),
    @("Dinamarca", 4077, 0, 1283, 2633, 142, 0, 161),
    @("Polonia", 3627, 0, 116, 3432, 50, 0, 79),
    @("Rumania", 3613, 0, 329, 3138, 119, 0, 146),
    @("India", 3588, 0, 229, 3260, 0, 0, 99),
    @("Malasia", 3483, 0, 915, 2511, 99, 0, 57),
    @("Ecuador", 3465, 0, 100, 3193, 100, 0, 172),
    @("Japon", 3139, 0, 514, 2548, 64, 0, 77),
    @("Filipinas", 3094, 0, 57, 2893, 1, 0, 144),
    @("Pakistan", 2818, 0, 131, 2646, 13, 0, 41),
    @("Luxemburgo", 2729, 0, 500, 2198, 35, 0, 31),
    @("Arabia Saudita", 2370, 191, 420, 1921, 41, 0, 29),
    @("Tailandia", 2169, 102, 674, 1472, 23, 3, 23),
    @("Indonesia", 2092, 0, 150, 1751, 0, 0, 191),
    @("Mexico", 1890, 202, 633, 1178, 1, 19, 79),
    @("Finlandia", 1882, 0, 300, 1557, 72, 0, 25),
    @("Panama", 1801, 0, 13, 1742, 50, 0, 46),
    @("Peru", 1746, 0, 914, 759, 88, 0, 73),
    @("Grecia", 1673, 0, 78, 1527, 92, 0, 68),
    @("Serbia", 1624, 0, 54, 1526, 89, 0, 44),
    @("Sudafrica", 1585, 0, 95, 1481, 7, 0, 9),
    @("Republica Dominicana", 1578, 0, 17, 1484, 147, 0, 77),
    @("Emiratos Arabes Unidos", 1505, 0, 125, 1370, 2, 0, 10),
    @("Argentina", 1451, 0, 279, 1129, 86, 0, 43),
    @("Islandia", 1417, 0, 396, 1017, 12, 0, 4),
    @("Colombia", 1406, 0, 85, 1289, 50, 0, 32),
    @("Catar", 1325, 0, 109, 1213, 37, 0, 3),
    @("Argelia", 1251, 0, 90, 1031, 0, 0, 130),
    @("Ucrania", 1225, 0, 25, 1168, 16, 0, 32),
    @("Singapur", 1189, 0, 297, 886, 24, 0, 6),
    @("Croacia", 1126, 0, 119, 995, 39, 0, 12),
    @("Egipto", 1070, 0, 241, 758, 0, 0, 71),
    @("Estonia", 1039, 0, 59, 967, 20, 0, 13),
    @("Nueva Zelanda", 1039, 89, 156, 882, 3, 0, 1),
    @("Eslovenia", 977, 0, 79, 876, 31, 0, 22),
    @("Marruecos", 919, 0, 66, 794, 1, 0, 59),
    @("Irak", 878, 0, 259, 563, 0, 0, 56),
    @("Hong Kong", 862, 0, 173, 685, 8, 0, 4),
    @("Lituania", 771, 0, 7, 753, 11, 0, 11),
    @("Armenia", 770, 0, 43, 720, 30, 0, 7),
    @("Moldavia", 752, 0, 29, 711, 65, 0, 12),
    @("Crucero", 712, 0, 619, 82, 10, 0, 11),
    @("Barein", 688, 0, 423, 261, 3, 0, 4),
    @("Hungria", 678, 0, 58, 588, 17, 0, 32),
    @("Bosnia y Herzegovina", 624, 0, 30, 573, 4, 0, 21),
    @("Camerun", 555, 0, 17, 529, 0, 0, 9),
    @("Tunez", 553, 0, 5, 530, 26, 0, 18),
    @("Kazajistan", 551, 20, 36, 510, 6, 0, 5),
    @("Azerbaiyan", 521, 0, 32, 484, 17, 0, 5),
    @("Libano", 520, 0, 54, 449, 27, 0, 17),
    @("Letonia", 509, 0, 1, 507, 3, 0, 1),
    @("Bulgaria", 503, 0, 34, 452, 26, 0, 17),
    @("Republica de Macedonia", 483, 0, 20, 446, 11, 0, 17),
    @("Kuwait", 479, 0, 93, 385, 17, 0, 1),
    @("Eslovaquia", 471, 0, 10, 460, 3, 0, 1),
    @("Principado de Andorra", 466, 0, 21, 428, 12, 0, 17),
    @("Bielorrusia", 440, 0, 53, 382, 11, 0, 5),
    @("Costa Rica", 435, 0, 13, 420, 13, 0, 2),
    @("Republica de Chipre", 426, 0, 33, 384, 11, 0, 9),
    @("Uruguay", 400, 0, 93, 302, 13, 0, 5),
    @("Taiwan", 355, 0, 50, 300, 0, 0, 5),
    @("Afganistan", 337, 28, 10, 320, 0, 0, 7),
    @("Reunion", 334, 0, 40, 294, 4, 0, 0),
    @("Albania", 333, 0, 99, 214, 7, 0, 20),
    @("Jordania", 323, 0, 74, 244, 5, 0, 5),
    @("Burkina Faso", 318, 0, 66, 236, 0, 0, 16),
    @("Cuba", 288, 0, 15, 267, 11, 0, 6),
    @("Oman", 277, 0, 61, 214, 3, 0, 2),
    @("Honduras", 268, 4, 6, 240, 10, 7, 22),
    @("Uzbekistan", 266, 0, 25, 239, 8, 0, 2),
    @("San Marino", 259, 0, 27, 200, 14, 0, 32),
    @("Costa de Marfil", 245, 0, 25, 219, 0, 0, 1),
    @("Vietnam", 240, 0, 90, 150, 3, 0, 0),
    @("Senegal", 219, 0, 72, 145, 1, 0, 2),
    @("Estado de Palestina", 217, 0, 21, 195, 0, 0, 1),
    @("Nigeria", 214, 0, 25, 185, 2, 0, 4),
    @("Malta", 213, 0, 2, 211, 3, 0, 0),
    @("Ghana", 205, 0, 31, 169, 2, 0, 5),
    @("Montenegro", 201, 0, 1, 198, 4, 0, 2),
    @("Mauricio", 196, 0, 7, 182, 1, 0, 7),
    @("Islas Feroe", 181, 0, 93, 88, 2, 0, 0),
    @("Sri Lanka", 166, 0, 27, 134, 5, 0, 5),
    @("Georgia", 162, 0, 36, 125, 6, 0, 1),
    @("Bolivia", 157, 18, 2, 145, 3, 0, 10),
    @("Venezuela", 155, 0, 52, 96, 6, 0, 7),
    @("Consejo Danes para los Refugiados", 154, 0, 3, 133, 0, 0, 18),
    @("Kirguistan", 147, 3, 9, 137, 5, 0, 1),
    @("Martinica", 145, 0, 27, 115, 22, 0, 3),
    @("Niger", 144, 0, 0, 136, 0, 0, 8),
    @("Brunei", 135, 0, 66, 68, 3, 0, 1),
    @("Mayotte", 134, 0, 14, 118, 3, 0, 2),
    @("Guadalupe", 134, 0, 24, 103, 14, 0, 7),
    @("Isla de Man", 126, 0, 0, 125, 0, 0, 1),
    @("Kenia", 126, 0, 4, 118, 2, 0, 4),
    @("Camboya", 114, 0, 50, 64, 1, 0, 0),
    @("Guinea", 111, 0, 5, 106, 0, 0, 0),
    @("Paraguay", 104, 8, 12, 89, 2, 0, 3),
    @("Trinidad yTobago", 103, 0, 1, 96, 0, 0, 6),
    @("Ruanda", 102, 0, 0, 102, 0, 0, 0),
    @("Gibraltar", 98, 0, 52, 46, 0, 0, 0),
    @("Liechtenstein", 77, 0, 0, 76, 0, 0, 1),
    @("Madagascar", 70, 0, 0, 70, 6, 0, 0),
    @("Banglades", 70, 0, 30, 32, 1, 0, 8),
    @("Monaco", 66, 0, 3, 62, 2, 0, 1),
    @("Aruba", 64, 0, 1, 63, 0, 0, 0),
    @("El Salvador", 62, 6, 2, 57, 4, 0, 3),
    @("Guatemala", 61, 0, 15, 44, 1, 0, 2),
    @("Guayana Francesa", 61, 0, 22, 39, 1, 0, 0),
    @("Jamaica", 55, 2, 7, 45, 0, 0, 3),
    @("Barbados", 52, 0, 0, 52, 0, 0, 0),
    @("Republica de Yibuti", 50, 0, 8, 42, 0, 0, 0),
    @("Uganda", 48, 0, 0, 48, 0, 0, 0),
    @("Macao", 44, 0, 10, 34, 0, 0, 0),
    @("Mali", 41, 0, 1, 37, 0, 0, 3),
    @("Togo", 41, 0, 17, 21, 0, 0, 3),
    @("Polinesia Francesa", 40, 0, 0, 40, 1, 0, 0),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Zambia", 39, 0, 2, 36, 0, 0, 1),
    @("Etiopia", 38, 0, 4, 34, 1, 0, 0),
    @("Bermudas", 37, 0, 14, 23, 0, 0, 0),
    @("Islas Caimanes", 35, 0, 1, 33, 0, 0, 1),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("Eritrea", 29, 0, 0, 29, 0, 0, 0),
    @("San Martin (Parte Francesa)", 29, 5, 7, 20, 6, 0, 2),
    @("Bahamas", 28, 0, 0, 24, 1, 0, 4),
    @("San Martin (Parte Holandesa)", 25, 2, 6, 15, 0, 0, 4),
    @("Guyana", 24, 1, 0, 20, 0, 0, 4),
    @("Congo", 22, 0, 2, 18, 0, 0, 2),
    @("Haiti", 21, 1, 1, 20, 0, 0, 0),
    @("Birmania", 21, 0, 0, 20, 0, 0, 1),
    @("Gabon", 21, 0, 1, 19, 0, 0, 1),
    @("Tanzania", 20, 0, 3, 16, 0, 0, 1),
    @("Maldivas", 19, 0, 13, 6, 0, 0, 0),
    @("Guinea-Bisau", 18, 0, 0, 18, 0, 0, 0),
    @("Libia", 18, 0, 0, 17, 0, 0, 1),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("Nueva Caledonia", 17, 0, 1, 16, 0, 0, 0),
    @("Guinea Ecuatorial", 16, 0, 1, 15, 0, 0, 0),
    @("Benin", 16, 0, 2, 14, 0, 0, 0),
    @("Siria", 16, 0, 2, 12, 0, 0, 2),
    @("Antigua y Barbuda", 15, 0, 0, 15, 1, 0, 0),
    @("Dominica", 14, 0, 0, 14, 0, 0, 0),
    @("Santa Lucia", 14, 0, 1, 13, 0, 0, 0),
    @("Mongolia", 14, 0, 2, 12, 0, 0, 0),
    @("Namibia", 14, 0, 3, 11, 0, 0, 0),
    @("Granada", 12, 0, 0, 12, 2, 0, 0),
    @("Fiyi", 12, 0, 0, 12, 0, 0, 0),
    @("Groenlandia", 11, 0, 3, 8, 0, 0, 0),
    @("Curazao", 11, 0, 5, 5, 0, 0, 1),
    @("Laos", 10, 0, 0, 10, 0, 0, 0),
    @("Seychelles", 10, 0, 0, 10, 0, 0, 0),
    @("Mozambique", 10, 0, 1, 9, 0, 0, 0),
    @("Surinam", 10, 0, 0, 9, 0, 0, 1),
    @("Sudan", 10, 0, 2, 6, 0, 0, 2),
    @("Liberia", 10, 0, 3, 6, 0, 0, 1),
    @("Angola", 10, 0, 2, 6, 0, 0, 2),
    @("Republica del Chad", 9, 0, 0, 9, 0, 0, 0),
    @("San Cristobal y Nieves", 9, 0, 0, 9, 0, 0, 0),
    @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
    @("Nepal", 9, 0, 1, 8, 0, 0, 0),
    @("Zimbabue", 9, 0, 0, 8, 0, 0, 1),
    @("Montserrat", 9, 0, 0, 7, 0, 0, 2),
    @("Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0),
    @("Santa Sede", 7, 0, 0, 7, 0, 0, 0),
    @("Somalia", 7, 0, 1, 6, 0, 0, 0),
    @("Cabo Verde", 7, 0, 0, 6, 0, 0, 1),
    @("San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0),
    @("San Bartolome", 6, 0, 1, 5, 0, 0, 0),
    @("Mauritania", 6, 0, 2, 3, 0, 0, 1),
    @("Islas Turcas y Caicos", 5, 0, 0, 5, 0, 0, 0),
    @("Nicaragua", 5, 0, 0, 4, 0, 0, 1),
    @("Butan", 5, 0, 2, 3, 0, 0, 0),
    @("Malaui", 4, 0, 0, 4, 0, 0, 0),
    @("Sahara Occidental", 4, 0, 0, 4, 0, 0, 0),
    @("Belice", 4, 0, 0, 4, 0, 0, 0),
    @("Sierra Leona", 4, 0, 0, 4, 0, 0, 0),
    @("Botsuana", 4, 0, 0, 3, 0, 0, 1),
    @("Gambia", 4, 0, 2, 1, 0, 0, 1),
    @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
    @("Burundi", 3, 0, 0, 3, 0, 0, 0),
    @("Anguila", 3, 0, 0, 3, 0, 0, 0),
    @("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0),
    @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
    @("Islas Malvinas", 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
